$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Edit 1 (slide 11): recolor "not all 19 SEE hold" from the literal
# purple RGB (7030A0) to the theme color accent6.
# -----------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(3)
$tr11 = $shp11.TextFrame.TextRange

$full11 = $tr11.Text
$target11 = "not all 19 SEE hold"
$idx11 = $full11.IndexOf($target11)
$run11 = $tr11.Characters($idx11 + 1, $target11.Length)
$run11.Font.Color.ObjectThemeColor = 10   # msoThemeColorAccent6

# -----------------------------------------------------------------
# Edit 2 (slide 3): collapse the 4 runs of
#   "With transition to " / "net zero" / ", some " / "SEE may fail to hold, thus:"
# into 2 runs:
#   "With transition to net zero" / ", some SEE may fail to hold, thus:"
# keeping the rPr of the 1st and 4th original runs respectively.
# -----------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(3)
$tr3 = $shp3.TextFrame.TextRange

$full3 = $tr3.Text
$runA = "With transition to "
$runB = "net zero"
$runC = ", some "
$runD = "SEE may fail to hold, thus:"

$startA = $full3.IndexOf($runA) + 1
$startB = $startA + $runA.Length
$startC = $startB + $runB.Length
$startD = $startC + $runC.Length

# Work from the end of the paragraph backwards so offsets computed
# above stay valid as each step is applied.
$rD = $tr3.Characters($startD, $runD.Length)
$rD.Text = ", some " + $runD

$rC = $tr3.Characters($startC, $runC.Length)
$rC.Text = ""

$rB = $tr3.Characters($startB, $runB.Length)
$rB.Text = ""

$rA = $tr3.Characters($startA, $runA.Length)
$rA.Text = "With transition to net zero"
